# The sheet's column A held a styled duplicate of column F (row numbers)
# that is no longer needed. Deleting it shifts B:F left into A:E, which
# keeps every other column (and its formatting/number values) intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:A").Delete()

# The header that used to read "MODEL_CONDITION" (now in column D after
# the shift) drops its underscore.
$ws.Range("D1").Value = "MODELCONDITION"

Write-Host "restructure complete"
